$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (codeforiati:group-code) and column D (codeforiati:group-name)
# had their contents swapped throughout the table, including the header row.
$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

$cVals = @{}
$dVals = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    $cVals[$r] = $ws.Cells.Item($r, 3).Value()
    $dVals[$r] = $ws.Cells.Item($r, 4).Value()
}
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $dVals[$r]
    $ws.Cells.Item($r, 4).Value = $cVals[$r]
}
